$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set E10, E12, E13, E14, E15 from FALSE to TRUE (boolean cells)
$ws.Range("E10").Value = $true
$ws.Range("E12").Value = $true
$ws.Range("E13").Value = $true
$ws.Range("E14").Value = $true
$ws.Range("E15").Value = $true
